$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column H (shifts H:L -> I:M)
$ws.Columns("H:H").Insert()

# Set the new header cell (copy style from the neighboring header, then set the text)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "barcode"

# Update the selection like the diff shows (I19)
$ws.Range("I19").Select()
